# Actualizo código del CRM
# Adds a new client record (C1002 / Tatiana Avila) to the "Clientes" sheet,
# as row 4, mirroring the existing row 3 (C1001 / Bandian Torres) record
# (same asesor, same dates, same estatus).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Seed row 4 by copying row 3 down first. This makes the new row inherit
# the same (unstyled) cell formatting/types as the other data rows, and
# in particular keeps the date-looking text values ("2025-10-08") stored
# as plain text instead of being auto-converted to Excel date serials,
# since fecha_ingreso/fecha_dispersion (columns E/F) end up identical to
# row 3's values and therefore don't need to be re-typed.
$ws.Range("A3:G3").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$ws.Application.CutCopyMode = $false

# Fill in the new client's distinguishing data.
$ws.Range("A4").Value = "C1002"
$ws.Range("B4").Value = "Tatiana Avila"
$ws.Range("C4").Value = "mundo E"
$ws.Range("D4").Value = "Martha Ortiz"
$ws.Range("G4").Value = "DISPERSADO"
